$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 - "Bom dia" / Ar da sala / A/C / 23 / FALSE / (blank)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Bom dia"
$ws.Range("B1").Value = "Ar da sala"
$ws.Range("C1").Value = "A/C"
$ws.Range("D1").Value = 23
$ws.Range("E1").Value = $false
$ws.Range("F1").WrapText = $false

# ---------------------------------------------------------------------
# Row 2 - "Bom dia" / tv da sala / Televisor / 1 / 100 / "false"
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Bom dia"
$ws.Range("B2").Value = "tv da sala"
$ws.Range("C2").Value = "Televisor"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = "'false"

# ---------------------------------------------------------------------
# Row 3 - "Boa noite" / tv da sala / Televisor / 1 / 0 / "false"
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Boa noite"
$ws.Range("B3").Value = "tv da sala"
$ws.Range("C3").Value = "Televisor"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "'false"

# ---------------------------------------------------------------------
# Row 4 - "Hora de estudar" / tv da sala / Televisor / 1 / 0 / FALSE
#         (kept on the default, un-formatted "Normal" style)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Hora de estudar"
$ws.Range("B4").Value = "tv da sala"
$ws.Range("C4").Value = "Televisor"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = $false
$ws.Range("A4:F4").Style = "Normal"

# ---------------------------------------------------------------------
# Number formatting: the two numeric columns (D, E) for the first three
# rows get a thousands-separator format, right aligned.
# ---------------------------------------------------------------------
$numRange = $ws.Range("D1:E3")
$numRange.NumberFormat = "#,##0"
$numRange.HorizontalAlignment = -4152

# A subset of those cells carries an explicit black font colour ...
$ws.Range("D1:D2").Font.Color = 0
$ws.Range("E2").Font.Color = 0

# ... while the remaining ones carry the (re-applied) theme font colour.
$ws.Range("E1").Font.Name = "Calibri"
$ws.Range("D3:E3").Font.Name = "Calibri"

# ---------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------
$ws.Rows(1).RowHeight = 19.5
$ws.Rows(2).RowHeight = 19.5
$ws.Rows(3).RowHeight = 18.75
